$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 40.88674644628622
$ws.Range("C8").Value = 50.27873197524158
$ws.Range("D8").Value = 98.82846243492406
$ws.Range("E8").Value = 98.91750232754872
$ws.Range("F8").Value = 98.47246877522802
$ws.Range("G8").Value = 98.02927001779099
$ws.Range("H8").Value = 97.56040581959799
$ws.Range("I8").Value = 96.28431623854225
